# Apply Mimecast email-report data refresh to the IT Metric Dashboard workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: Meraki Switches ---
$wsSwitches = $wb.Worksheets.Item("Meraki Switches")
$wsSwitches.Range("C21").Value = 100
$wsSwitches.Range("A30").Value = "DEN-LAB"

# --- Sheet: Meraki AP ---
$wsAP = $wb.Worksheets.Item("Meraki AP")

$wsAP.Range("C2").Value = 62.3
$wsAP.Range("D2").Value = "online"

$wsAP.Range("C3").Value = 62.02
$wsAP.Range("D3").Value = "online"

$wsAP.Range("C11").Value = 62.25
$wsAP.Range("D11").Value = "online"

$wsAP.Range("C12").Value = 62.3
$wsAP.Range("D12").Value = "online"

$wsAP.Range("C14").Value = 62.13
$wsAP.Range("D14").Value = "online"

$wsAP.Range("C15").Value = 62.28
$wsAP.Range("D15").Value = "online"

$wsAP.Range("C16").Value = 62.3
$wsAP.Range("D16").Value = "online"

$wsAP.Range("C17").Value = 62.3
$wsAP.Range("D17").Value = "online"

$wsAP.Range("C18").Value = 62.3
$wsAP.Range("D18").Value = "online"

$wsAP.Range("C19").Value = 62.29
$wsAP.Range("D19").Value = "online"

$wsAP.Range("C20").Value = 62.3
$wsAP.Range("D20").Value = "online"

$wsAP.Range("C21").Value = 62.3
$wsAP.Range("D21").Value = "online"

$wsAP.Range("C22").Value = 62.3
$wsAP.Range("D22").Value = "online"

$wsAP.Range("C23").Value = 62.29
$wsAP.Range("D23").Value = "online"

$wsAP.Range("C24").Value = 62.3
$wsAP.Range("D24").Value = "online"

$wsAP.Range("C25").Value = 62.29
$wsAP.Range("D25").Value = "online"

$wsAP.Range("C26").Value = 62.29
$wsAP.Range("D26").Value = "online"

# --- Sheet: Freshservice ---
$wsFS = $wb.Worksheets.Item("Freshservice")
$wsFS.Range("D1").Value = "Resolution Rate"
$wsFS.Range("A2").Value = 51
$wsFS.Range("B2").Value = 10
$wsFS.Range("C2").Value = 41
$wsFS.Range("D2").Value = 80.39
